$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.293.45"
$ws.Range("E2").Value = "  -0.07%  "

# Row 3
$ws.Range("D3").Value = "3.682.51"

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "682.56"
$ws.Range("E5").Value = "  -0.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.06"
$ws.Range("E6").Value = "  -2.74%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  -1.16%  "

# Row 9
$ws.Range("E9").Value = "  -1.54%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.01"
$ws.Range("E10").Value = "  -2.77%  "

# Row 11
$ws.Range("E11").Value = "  -2.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000232"
$ws.Range("E12").Value = "  -1.94%  "

# Row 13
$ws.Range("D13").Value = "4.303.50"
$ws.Range("E13").Value = "  +0.01%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.27"
$ws.Range("E14").Value = "  -3.75%  "

# Row 15
$ws.Range("D15").Value = "3.720.90"
$ws.Range("E15").Value = "  +0.98%  "

# Row 16
$ws.Range("D16").Value = "69.311.35"
$ws.Range("E16").Value = "  -0.07%  "

# Row 17
$ws.Range("E17").Value = "  +1.77%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.95"
$ws.Range("E18").Value = "  -2.39%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.38"
$ws.Range("E19").Value = "  -4.16%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.83"
$ws.Range("E20").Value = "  -1.66%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.649"
$ws.Range("E22").Value = "  -2.74%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.99"
$ws.Range("E23").Value = "  -0.18%  "

# Row 24
$ws.Range("D24").Value = "3.827.07"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("E25").Value = "  +0.00%  "

# Row 26
$ws.Range("E26").Value = "  -5.93%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.90"
$ws.Range("E27").Value = "  -5.39%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.14"
$ws.Range("E28").Value = "  -4.67%  "

# Row 29
$ws.Range("E29").Value = "  -1.65%  "

# Row 30
$ws.Range("E30").Value = "  -5.15%  "

# Row 31
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.56"
$ws.Range("E31").Value = "  -4.71%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.98"
$ws.Range("E32").Value = "  -6.01%  "

# Row 33
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.995"
$ws.Range("E33").Value = "  -0.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.84"
$ws.Range("E34").Value = "  -1.10%  "

# Row 35
$ws.Range("D35").Value = "3.662.06"
$ws.Range("E35").Value = "  +0.52%  "

# Row 36
$ws.Range("E36").Value = "  -3.79%  "

# Row 37
$ws.Range("E37").Value = "  -4.55%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.05"
$ws.Range("E38").Value = "  -1.34%  "

# Row 40
$ws.Range("E40").Value = "  +3.46%  "

# Row 41
$ws.Range("E41").Value = "  -4.00%  "

# Row 42
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "168.87"
$ws.Range("E43").Value = "  +9.57%  "

# Row 44
$ws.Range("E44").Value = "  -1.83%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.56"
$ws.Range("E45").Value = "  -1.36%  "

# Row 46
$ws.Range("E46").Value = "  -4.98%  "

# Row 47
$ws.Range("E47").Value = "  -1.39%  "

# Row 48
$ws.Range("E48").Value = "  +2.56%  "

# Row 49
$ws.Range("E49").Value = "  -4.77%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.76"
$ws.Range("E50").Value = "  -4.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.266"
$ws.Range("E51").Value = "  -2.22%  "
